$d = $word.ActiveDocument

# 1) "Switch for detecting box in effector." -> "Switch for effector. Active LOW."
#    (unique text in the document, safe to replace across whole story)
$d.Content.Find.Execute("for detecting box in effector.", $true, $false, $false, $false, $false,
                         $true, 1, $false, "for effector. Active LOW.", 2)

# 2) Expansion Hub / twobar row: Port "Motor 1" -> "Motor 3"
#    ("Motor 1" also appears in the frMotor row, so the replace must be
#    confined to the twobar row's Port cell; wdReplaceAll searches the
#    whole story regardless of range bounds, so use wdReplaceOne on a
#    range rebuilt from the cell's Start/End which IS bounds-respecting)
$t = $d.Tables.Item(1)
for ($i = 1; $i -le $t.Rows.Count; $i++) {
    $row = $t.Rows.Item($i)
    $nameCell = $row.Cells.Item(2).Range.Text
    if ($nameCell -like "*twobar*") {
        $cellRange = $row.Cells.Item(3).Range
        $rng = $d.Range($cellRange.Start, $cellRange.End)
        $rng.Find.Execute("Motor 1", $true, $false, $false, $false, $false,
                           $true, 1, $false, "Motor 3", 1)
    }
}
